$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update hashcode values (column B) for the rows whose code (column A) matches.
# Each entry: row number, expected code in column A, new hashcode for column B.
$updates = @(
    @{ Row = 34; Code = "05-050316TP"; NewHash = "c61e0c5fa0c3d3aeb7f195c62229f494" }
    @{ Row = 44; Code = "05-050105A"; NewHash = "a2cfcbfef9b7b4aed5ed06cdf76e820f" }
    @{ Row = 74; Code = "05-050103A"; NewHash = "9555bf74da8a390313ded720eb47dce7" }
    @{ Row = 89; Code = "05-050104A"; NewHash = "160ee88f449d69ffbf488ebe9d2dcc44" }
    @{ Row = 99; Code = "05-050101A"; NewHash = "ec5bd2a050b8a245967e920be6cdaaa2" }
    @{ Row = 110; Code = "05-050102A"; NewHash = "4050bd447a74401c61ea746f9711d4fc" }
    @{ Row = 154; Code = "05-050007TC"; NewHash = "e9828e955ed4896624069e2230da5da2" }
    @{ Row = 160; Code = "05-050007TP"; NewHash = "f3de5288eeaf606f566c40f38f1f948a" }
    @{ Row = 161; Code = "05-050105TC"; NewHash = "9bb4c7968671c6ffbee5b3db18131f17" }
    @{ Row = 162; Code = "05-050308A"; NewHash = "28b7081ddd8b2bf574091a34d8703cef" }
    @{ Row = 168; Code = "05-050105TP"; NewHash = "36c8cd53ba8a46717318adc0a51706b1" }
    @{ Row = 180; Code = "05-050303TC"; NewHash = "4452182d4a3e39871668d09fdb6c1e5b" }
    @{ Row = 213; Code = "05-050303A"; NewHash = "e11742ebab986b101aaf472dd8371e81" }
    @{ Row = 278; Code = "05-050101TP"; NewHash = "4f4e6e1d7f91885a3a4f184b8ac396e3" }
    @{ Row = 330; Code = "05-050005TC"; NewHash = "0f541db1bee54323ba63ade30ce40dfc" }
    @{ Row = 335; Code = "05-050005TP"; NewHash = "ecbe729ac86df7acbe5e7934836f2f14" }
    @{ Row = 345; Code = "05-050103TP"; NewHash = "183913fecc02620ae6913e0667b17656" }
    @{ Row = 461; Code = "05-050313A"; NewHash = "b11b80ec3b93464d6b97a5f9c1948435" }
    @{ Row = 506; Code = "05-050306TP"; NewHash = "51d94fbb108c060af0774f3dfc25fd2e" }
    @{ Row = 514; Code = "05-050317TC"; NewHash = "1ff4dd27e25e4cecffa8c888a063c5c2" }
    @{ Row = 524; Code = "05-050317TP"; NewHash = "586802b4d9ba45de50d961c63708f3c0" }
    @{ Row = 534; Code = "05-050006A"; NewHash = "76da3783aa2a61aa6867b6ba825b3179" }
    @{ Row = 547; Code = "05-050201A"; NewHash = "12134a6651c6de21c72dc6c1e1dae89a" }
    @{ Row = 553; Code = "05-050007A"; NewHash = "58d85ba2051dd71507a5e4255d2e5b94" }
    @{ Row = 584; Code = "05-050005A"; NewHash = "a576e1b2662d1a21d6c1d37626fd4452" }
    @{ Row = 666; Code = "05-050317A"; NewHash = "6a504f8d367e29df8fe91b6e061f2350" }
    @{ Row = 729; Code = "05-050316A"; NewHash = "27ed38bf1fbffac7273df8279ccba7ca" }
    @{ Row = 768; Code = "05-050102TP"; NewHash = "8a866f38cea4d509d812189b47eef642" }
    @{ Row = 811; Code = "05-050006TC"; NewHash = "dbd952bba9bedbb15ced3d14a76bc9b0" }
    @{ Row = 815; Code = "05-050006TP"; NewHash = "bd5b9380588c9dc7c9ba8123dc3cab76" }
    @{ Row = 816; Code = "05-050104TC"; NewHash = "1951623ae9020a139ec3467817acc2ab" }
    @{ Row = 825; Code = "05-050104TM"; NewHash = "76fb08e3968f1341beee8c4d704ab1a6" }
    @{ Row = 827; Code = "05-050104TP"; NewHash = "fe391b223dd9b3e7fc6a5f6ebd9890a3" }
)

foreach ($u in $updates) {
    $actualCode = $ws.Cells.Item($u.Row, 1).Value2
    if ($actualCode -ne $u.Code) {
        Write-Host "WARNING: row $($u.Row) expected code '$($u.Code)' but found '$actualCode'"
    }
    $ws.Cells.Item($u.Row, 2).Value = $u.NewHash
}
